$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 292, pushing the existing rows 292-299
# down to 295-302 (formatting of row 292 - e.g. the date style on column D -
# is inherited by the newly inserted rows).
$ws.Rows("292:294").Insert()

# --- Row 292 (new) ---
$ws.Range("A292").Value = 4
$ws.Range("B292").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C292").Value = "Los Lagos"
$ws.Range("D292").Value = 44939
$ws.Range("E292").Value = 10
$ws.Range("F292").Value = 100112024
$ws.Range("G292").Value = "Choclo"
$ws.Range("H292").Value = "Choclero"
$ws.Range("I292").Value = "Primera"
$ws.Range("J292").Value = 9000
$ws.Range("K292").Value = 500
$ws.Range("L292").Value = 500
$ws.Range("M292").Value = 500
$ws.Range("N292").Value = "$/unidad"
$ws.Range("O292").Value = "Región de O'Higgins"
$ws.Range("P292").Value = 500
$ws.Range("Q292").Value = 1
$ws.Range("R292").Value = "Hortaliza"

# --- Row 293 (new) ---
$ws.Range("A293").Value = 4
$ws.Range("B293").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C293").Value = "Los Lagos"
$ws.Range("D293").Value = 44939
$ws.Range("E293").Value = 10
$ws.Range("F293").Value = 100112024
$ws.Range("G293").Value = "Choclo"
$ws.Range("H293").Value = "Choclero"
$ws.Range("I293").Value = "Segunda"
$ws.Range("J293").Value = 9000
$ws.Range("K293").Value = 350
$ws.Range("L293").Value = 350
$ws.Range("M293").Value = 350
$ws.Range("N293").Value = "$/unidad"
$ws.Range("O293").Value = "Región de O'Higgins"
$ws.Range("P293").Value = 350
$ws.Range("Q293").Value = 1
$ws.Range("R293").Value = "Hortaliza"

# --- Row 294 (new) ---
$ws.Range("A294").Value = 4
$ws.Range("B294").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C294").Value = "Los Lagos"
$ws.Range("D294").Value = 44939
$ws.Range("E294").Value = 10
$ws.Range("F294").Value = 100112024
$ws.Range("G294").Value = "Choclo"
$ws.Range("H294").Value = "Dulce o Americano"
$ws.Range("I294").Value = "Primera"
$ws.Range("J294").Value = 25000
$ws.Range("K294").Value = 250
$ws.Range("L294").Value = 250
$ws.Range("M294").Value = 250
$ws.Range("N294").Value = "$/unidad"
$ws.Range("O294").Value = "Región del Maule"
$ws.Range("P294").Value = 250
$ws.Range("Q294").Value = 1
$ws.Range("R294").Value = "Hortaliza"
